$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "January 1-June 30, 2023"

$data = @{
    3 = @(34460, 5380, 5666)
    4 = @(20702, 1990, 2449)
    5 = @(56847, 5220, 5198)
    6 = @(1146, 554, 150)
    7 = @(38343, 6383, 4650)
    8 = @(4588, 805, 1056)
    9 = @(4509, 968, 625)
    10 = @(2403, 316, 173)
    11 = @(701, 247, 2)
    12 = @(0, 0, 0)
    13 = @(797, 154, 244)
    14 = @(2393, 867, 934)
    15 = @(3810, 1366, 715)
    16 = @(2968, 1300, 345)
    17 = @(2279, 638, 277)
    18 = @(14054, 1984, 2604)
    19 = @(1087, 518, 274)
    20 = @(14172, 2027, 2280)
    21 = @(173, 338, 13)
    22 = @(13671, 1629, 2241)
    23 = @(858, 544, 133)
    24 = @(15137, 1626, 2934)
    25 = @(62020, 5548, 7176)
    26 = @(5031, 1544, 673)
    27 = @(0, 0, 0)
    28 = @(4054, 958, 1002)
    29 = @(1576, 451, 332)
    30 = @(11198, 1954, 2222)
    31 = @(407, 93, 244)
    32 = @(2557, 1363, 233)
    33 = @(11070, 2617, 2318)
    34 = @(8589, 2645, 1545)
    35 = @(4014, 445, 766)
    36 = @(46193, 4937, 4094)
    37 = @(6562, 2224, 889)
    38 = @(21548, 1396, 2015)
    39 = @(802, 847, 148)
    40 = @(1023, 416, 404)
    41 = @(1938, 462, 68)
    42 = @(7403, 353, 268)
    43 = @(177, 158, 13)
    44 = @(514, 114, 49)
    45 = @(0, 0, 0)
    46 = @(2664, 890, 339)
    47 = @(10714, 2686, 2006)
    48 = @(27976, 2623, 4083)
    49 = @(12624, 2634, 1214)
    50 = @(9275, 1192, 1985)
    51 = @(27108, 2450, 4271)
    52 = @(3842, 829, 930)
    53 = @(9113, 1745, 1441)
    54 = @(1761, 1085, 632)
    55 = @(1844, 1106, 98)
    56 = @(4041, 947, 1459)
    57 = @(9968, 3687, 2238)
    58 = @(10032, 1087, 398)
    59 = @(542736, 81566, 75880)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
